$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 data: 번호(A), 문제 이름(B), 틀린 날짜(C)
$ws.Range("A15").Value = 1987
$ws.Range("B15").Value = "알파벳(비트마스킹)"

# Copy the date cell's style from the row above (C14) so no new numFmt/style is created
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C15").Value = 45910

# Update the frozen pane / selection view to match the new data extent
$ws.Range("B16").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
